# Refresh each cryptocurrency row with the latest scraped Price (D) / Volume(1h) (E)
# text. Values are written through a formula-literal + paste-as-values round trip so
# that numeric-looking strings (e.g. "0.7062", "9.000") land back in the sheet as plain
# text (matching the scraper's inline-string output) instead of being auto-coerced to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @(2, 4, "=`"29.414.92`""),
    @(2, 5, "=`"  +0.26%  `""),
    @(3, 4, "=`"1.865.96`""),
    @(3, 5, "=`"  -0.66%  `""),
    @(5, 4, "=`"0.7062`""),
    @(5, 5, "=`"  -0.63%  `""),
    @(6, 4, "=`"243.05`""),
    @(6, 5, "=`"  +0.21%  `""),
    @(7, 5, "=`"  +0.00%  `""),
    @(8, 4, "=`"0.07876`""),
    @(8, 5, "=`"  -2.10%  `""),
    @(9, 4, "=`"0.3123`""),
    @(9, 5, "=`"  -1.59%  `""),
    @(10, 5, "=`"  -2.32%  `""),
    @(11, 4, "=`"0.08011`""),
    @(11, 5, "=`"  -3.71%  `""),
    @(12, 4, "=`"1.870.07`""),
    @(12, 5, "=`"  -0.63%  `""),
    @(13, 4, "=`"5.198`""),
    @(13, 5, "=`"  -1.19%  `""),
    @(14, 4, "=`"93.49`""),
    @(14, 5, "=`"  -1.29%  `""),
    @(15, 4, "=`"0.6987`""),
    @(15, 5, "=`"  -2.62%  `""),
    @(16, 4, "=`"6.447`""),
    @(16, 5, "=`"  +0.67%  `""),
    @(17, 4, "=`"0.000008366`""),
    @(17, 5, "=`"  -3.60%  `""),
    @(18, 4, "=`"29.398.05`""),
    @(18, 5, "=`"  +0.21%  `""),
    @(19, 4, "=`"253.09`""),
    @(19, 5, "=`"  +4.32%  `""),
    @(20, 4, "=`"2.122.12`""),
    @(20, 5, "=`"  -0.64%  `""),
    @(21, 4, "=`"13.10`""),
    @(21, 5, "=`"  -1.79%  `""),
    @(22, 5, "=`"  +0.02%  `""),
    @(23, 4, "=`"7.615`""),
    @(23, 5, "=`"  -2.63%  `""),
    @(24, 5, "=`"  +0.04%  `""),
    @(25, 4, "=`"0.1558`""),
    @(25, 5, "=`"  -0.75%  `""),
    @(26, 4, "=`"9.000`""),
    @(26, 5, "=`"  -1.08%  `""),
    @(27, 4, "=`"160.64`""),
    @(27, 5, "=`"  -1.52%  `""),
    @(28, 4, "=`"18.74`""),
    @(28, 5, "=`"  +0.91%  `""),
    @(29, 4, "=`"1.499`""),
    @(29, 5, "=`"  -0.75%  `""),
    @(30, 4, "=`"4.321`""),
    @(30, 5, "=`"  -2.53%  `""),
    @(31, 4, "=`"4.283`""),
    @(31, 5, "=`"  -1.31%  `""),
    @(32, 4, "=`"1.209`""),
    @(32, 5, "=`"  +1.04%  `""),
    @(33, 4, "=`"0.05308`""),
    @(33, 5, "=`"  -1.88%  `""),
    @(34, 4, "=`"1.887`""),
    @(34, 5, "=`"  -2.98%  `""),
    @(35, 4, "=`"0.7511`""),
    @(35, 5, "=`"  -2.79%  `""),
    @(36, 4, "=`"1.167`""),
    @(36, 5, "=`"  -1.78%  `""),
    @(37, 4, "=`"2.710`""),
    @(37, 5, "=`"  +0.96%  `""),
    @(38, 4, "=`"0.01876`""),
    @(38, 5, "=`"  -0.73%  `""),
    @(39, 4, "=`"1.274.82`""),
    @(39, 5, "=`"  +0.72%  `""),
    @(40, 4, "=`"2.741`""),
    @(40, 5, "=`"  -0.32%  `""),
    @(41, 4, "=`"0.8969`""),
    @(41, 5, "=`"  -1.18%  `""),
    @(42, 4, "=`"108.84`""),
    @(42, 5, "=`"  -4.08%  `""),
    @(43, 4, "=`"6.003`""),
    @(43, 5, "=`"  -7.56%  `""),
    @(44, 4, "=`"71.19`""),
    @(44, 5, "=`"  -4.64%  `""),
    @(45, 5, "=`"  +0.01%  `""),
    @(46, 5, "=`"  -3.03%  `""),
    @(47, 4, "=`"2.021.59`""),
    @(47, 5, "=`"  -0.26%  `""),
    @(48, 4, "=`"9.588`""),
    @(48, 5, "=`"  +0.66%  `""),
    @(49, 4, "=`"1.787`""),
    @(49, 5, "=`"  -1.32%  `""),
    @(50, 4, "=`"0.5168`""),
    @(50, 5, "=`"  -1.03%  `""),
    @(51, 4, "=`"0.4304`""),
    @(51, 5, "=`"  -1.57%  `""),
)

foreach ($u in $updates) {
    $c = $ws.Cells.Item($u[0], $u[1])
    $c.Formula = $u[2]
    $c.Copy()
    $c.PasteSpecial(-4163)  # xlPasteValues
}

$excel.CutCopyMode = $false
